# Updated cryptos list on Tue Apr  4 23:40:00 UTC 2023 with GitHub Actions
# Refresh the Price (D) and Volume(1h) (E) columns for the cryptos table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.182.19"
$ws.Range("E2").Value = "  +1.21%  "
$ws.Range("D3").Value = "1.871.02"
$ws.Range("E3").Value = "  +3.35%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'311.97"
$ws.Range("E5").Value = "  +0.95%  "
$ws.Range("D7").Value = "'0.5043"
$ws.Range("E7").Value = "  +1.02%  "
$ws.Range("D8").Value = "'0.3898"
$ws.Range("E8").Value = "  +0.29%  "
$ws.Range("D9").Value = "'0.09498"
$ws.Range("E9").Value = "  -0.81%  "
$ws.Range("D10").Value = "'1.141"
$ws.Range("E10").Value = "  +3.83%  "
$ws.Range("D11").Value = "'40.86"
$ws.Range("E11").Value = "  +1.33%  "
$ws.Range("D12").Value = "'6.442"
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("D13").Value = "'20.91"
$ws.Range("E13").Value = "  +2.23%  "
$ws.Range("D14").Value = "1.873.05"
$ws.Range("E14").Value = "  +3.16%  "
$ws.Range("D15").Value = "'1.003"
$ws.Range("E15").Value = "  +0.14%  "
$ws.Range("D16").Value = "'7.378"
$ws.Range("E16").Value = "  +1.71%  "
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("D18").Value = "'92.66"
$ws.Range("E18").Value = "  -0.75%  "
$ws.Range("D19").Value = "'0.06605"
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("E20").Value = "  +3.10%  "
$ws.Range("D21").Value = "'1.002"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").Value = "'6.176"
$ws.Range("E22").Value = "  +3.90%  "
$ws.Range("D23").Value = "28.250.38"
$ws.Range("E23").Value = "  +1.20%  "
$ws.Range("D24").Value = "'11.25"
$ws.Range("E24").Value = "  +0.78%  "
$ws.Range("D25").Value = "'2.279"
$ws.Range("E25").Value = "  +1.37%  "
$ws.Range("D26").Value = "'2.576"
$ws.Range("E26").Value = "  +7.33%  "
$ws.Range("D27").Value = "2.086.67"
$ws.Range("E27").Value = "  +3.14%  "
$ws.Range("E28").Value = "  +1.87%  "
$ws.Range("D29").Value = "'159.13"
$ws.Range("E29").Value = "  +0.95%  "
$ws.Range("D30").Value = "'127.17"
$ws.Range("E30").Value = "  -0.83%  "
$ws.Range("D31").Value = "'0.1061"
$ws.Range("E31").Value = "  -0.94%  "
$ws.Range("D32").Value = "'1.065"
$ws.Range("E32").Value = "  +1.20%  "
$ws.Range("E33").Value = "  +0.57%  "
$ws.Range("D35").Value = "'0.06739"
$ws.Range("E35").Value = "  -0.99%  "
$ws.Range("D36").Value = "'9.509"
$ws.Range("E36").Value = "  +5.60%  "
$ws.Range("D37").Value = "'0.02407"
$ws.Range("E37").Value = "  +3.87%  "
$ws.Range("D38").Value = "'0.2184"
$ws.Range("E38").Value = "  +1.85%  "
$ws.Range("E39").Value = "  +1.18%  "
$ws.Range("D40").Value = "'0.6352"
$ws.Range("E40").Value = "  +1.77%  "
$ws.Range("D41").Value = "'4.996"
$ws.Range("E41").Value = "  +1.27%  "
$ws.Range("D42").Value = "'1.186"
$ws.Range("E42").Value = "  +3.58%  "
$ws.Range("D43").Value = "'1.002"
$ws.Range("E43").Value = "  +0.20%  "
$ws.Range("D44").Value = "'13.46"
$ws.Range("E44").Value = "  +2.98%  "
$ws.Range("D45").Value = "'0.5979"
$ws.Range("E45").Value = "  +1.08%  "
$ws.Range("E46").Value = "  -1.20%  "
$ws.Range("E47").Value = "  -0.73%  "
$ws.Range("D48").Value = "'1.994"
$ws.Range("E48").Value = "  +1.90%  "
$ws.Range("D49").Value = "'123.39"
$ws.Range("E49").Value = "  -0.35%  "
$ws.Range("E50").Value = "  +1.59%  "
$ws.Range("D51").Value = "'0.06849"
